$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values
# (e.g. "560.98") are not auto-converted to floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.894.01"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "560.98"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").Value = "142.33"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("E11").Value = "  -2.70%  "
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("D15").Value = "2.832.47"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").Value = "61.795.17"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "2.398.68"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "11.18"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").Value = "320.33"
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").Value = "6.78"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").Value = "65.80"
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("E24").Value = "  -2.86%  "
$ws.Range("D25").Value = "8.77"
$ws.Range("E25").Value = "  -4.40%  "
$ws.Range("D26").Value = "560.45"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D28").Value = "2.520.13"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").Value = "0.0₃0929"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").Value = "8.15"
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("E31").Value = "  -4.78%  "
$ws.Range("D32").Value = "0.146"
$ws.Range("E32").Value = "  -1.71%  "
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  -4.86%  "
$ws.Range("E36").Value = "  -2.44%  "
$ws.Range("D37").Value = "152.37"
$ws.Range("E37").Value = "  +2.87%  "
$ws.Range("E38").Value = "  -5.76%  "
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").Value = "18.51"
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("E41").Value = "  -5.95%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  -3.35%  "
$ws.Range("D44").Value = "147.02"
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("D47").Value = "19.72"
$ws.Range("E47").Value = "  -3.57%  "
$ws.Range("D48").Value = "0.585"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("D49").Value = "0.0915"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("E51").Value = "  +0.46%  "

# Restore the original (default) style for column D now that the
# text values are safely stored, avoiding a lingering number format.
$ws.Range("D2:D51").Style = "Normal"
